# TimberBeamData.xlsx: add a couple of helper columns (B, C) next to the
# existing FirstName/value column A, then total/average them in row 4.
# The SUM cell (B4) is the value that gets passed into the PDF export, so
# it is colored red to flag it; C4 picks up an explicit font too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data values in columns B and C (rows 2-3)
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 4
$ws.Range("B3").Value = 3
$ws.Range("C3").Value = 6

# Row 4: totals for the new columns
$ws.Range("B4").Formula = "=SUM(B2:B3)"
$ws.Range("C4").Formula = "=AVERAGE(C2:C3)"

# B4 is the cell whose value is read for the PDF - make it red so it's easy
# to spot in the sheet.
$ws.Range("B4").Font.Color = 255
# C4 also gets an explicit font applied.
$ws.Range("C4").Font.Name = "Calibri"

# Cursor ends up on E12 after the edits
$ws.Range("E12").Select()

# Page setup touched (orientation recorded) as part of preparing the sheet
# for printing/export to PDF.
$ws.PageSetup.Orientation = 1
$ws.PageSetup.HorizontalDpi = 300
$ws.PageSetup.VerticalDpi = 300

$wb.Save()
